$d = $word.ActiveDocument

# The checklist header table is the first table in the document:
#   Row 2, Col 4 -> "Sprint No."  value  (currently "1")
#   Row 3, Col 2 -> "Review Date" value  (currently "02/09/18", merged across cols 2-4)
$t = $d.Tables.Item(1)

# --- Update Sprint No. from "1" to "2" --------------------------------
$sprintCell = $t.Cell(2, 4)
$sprintRange = $sprintCell.Range
$sprintRange.Find.Execute("1", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "2", 1) | Out-Null

# --- Update Review Date from "02/09/18" to "02/21/18" ------------------
$t2 = $d.Tables.Item(1)
$dateCell = $t2.Cell(3, 2)
$dateRange = $dateCell.Range
$dateRange.Find.Execute("02/09/18", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "02/21/18", 1) | Out-Null
